# Generate Report for Handoff
# File "b.md" has completed a new handoff round (row 3 on each sheet).
# Update status + the new handoff file name / timestamp for both locales,
# and reflect the same "Ready for handoff" status + date on the Overview sheet.

$wb = $excel.ActiveWorkbook

# NOTE: to update a hyperlinked cell's displayed text AND keep the
# hyperlink's "display" attribute (and r:id / formatting) in sync without
# creating a duplicate hyperlink entry, the Hyperlink object must be
# captured from the worksheet's `Hyperlinks` collection via `foreach`;
# going through `.Hyperlinks.Item(n)` (or `Range(...).Hyperlinks.Item(n)`)
# does not resolve properties correctly in this host.

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is "b.md"
#   B3 (zh-cn status)  : Handed back... -> Ready for handoff
#   C3 (de-de status)  : Handed back... -> Ready for handoff
#   D3 (Latest Handoff Date): -> 2016-34-13 06:34:49
# (none of these three cells are hyperlinked)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-34-13 06:34:49"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is "b.md"
#   C3 (Status)               : Handed back... -> Ready for handoff
#   D3 (Latest Handoff File)  : a....zh-cn.xlf -> b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf (hyperlinked)
#   E3 (Latest Handoff Datetime): -> 2016-03-13 06:34:45
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

$wsZhCn.Range("E3").Value = "2016-03-13 06:34:45"

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is "b.md"
#   C3 (Status)               : Handed back... -> Ready for handoff
#   D3 (Latest Handoff File)  : a....de-de.xlf -> b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf (hyperlinked)
#   E3 (Latest Handoff Datetime): -> 2016-03-13 06:34:49
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}

$wsDeDe.Range("E3").Value = "2016-03-13 06:34:49"
